$wb = $excel.ActiveWorkbook

# --- "RR death by anemia" sheet: add population-group columns -------------
# Currently the sheet only has age-band columns (C:H = <1 month .. 40-44
# years headers mixed together). Insert 6 new columns before the existing
# age-band columns so the sheet gets the full set of population groups:
#   child age bands, pregnant women, then women-of-reproductive-age bands.
$wsRR = $wb.Worksheets.Item("RR death by anemia")

$wsRR.Range("C1:H1").EntireColumn.Insert()

$wsRR.Cells.Item(1, 3).Value2 = "<1 month"
$wsRR.Cells.Item(1, 4).Value2 = "1-5 months"
$wsRR.Cells.Item(1, 5).Value2 = "6-11 months"
$wsRR.Cells.Item(1, 6).Value2 = "12-23 months"
$wsRR.Cells.Item(1, 7).Value2 = "24-59 months"
$wsRR.Cells.Item(1, 8).Value2 = "pregnant women"

for ($r = 2; $r -le 7; $r++) {
    for ($c = 3; $c -le 8; $c++) {
        $wsRR.Cells.Item($r, $c).Value2 = 1
    }
}

# --- "anemia prevalence" sheet: selection moves from J24 to G15 -----------
$wsAnemia = $wb.Worksheets.Item("anemia prevalence")
$wsAnemia.Activate()
$wsAnemia.Range("G15").Select()

# --- Make "RR death by anemia" the active tab/sheet (also restores its
#     own selection at F20 and drops the stale tabSelected flag that was
#     previously sitting on "distributions") ---------------------------
$wsRR.Activate()
$wsRR.Range("F20").Select()
